$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.597
$ws.Range("C2").Value = 1.609
$ws.Range("D2").Value = 1.191
$ws.Range("E2").Value = 2.191
$ws.Range("F2").Value = 1.826

$ws.Range("B3").Value = 3.968
$ws.Range("C3").Value = 4.066
$ws.Range("D3").Value = 2.799
$ws.Range("E3").Value = 9.265000000000001
$ws.Range("F3").Value = 5.964

$ws.Range("B4").Value = 18.214
$ws.Range("C4").Value = 15.308
$ws.Range("D4").Value = 12.948
$ws.Range("E4").Value = 30.588
$ws.Range("F4").Value = 18.064

$ws.Range("B5").Value = 17.792
$ws.Range("C5").Value = 15.51
$ws.Range("D5").Value = 12.991
$ws.Range("E5").Value = 29.278
$ws.Range("F5").Value = 17.525

$ws.Range("B6").Value = 0.663
$ws.Range("C6").Value = 0.76
$ws.Range("D6").Value = 0.483
$ws.Range("E6").Value = 0.786
$ws.Range("F6").Value = 0.732
